$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 650-651, pushing the existing rows 650-686 down to 652-688.
# EntireRow.Insert() carries formatting (incl. the date style on column D) to the
# freshly inserted rows, matching the surrounding data rows.
$ws.Range("A650:R651").EntireRow.Insert()

# --- New row 650: Lechuga / Conconina(o) ---
$ws.Range("A650").Value = 7
$ws.Range("B650").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C650").Value = "Ñuble"
$ws.Range("D650").Value = 44753
$ws.Range("E650").Value = 16
$ws.Range("F650").Value = 100112033
$ws.Range("G650").Value = "Lechuga"
$ws.Range("H650").Value = "Conconina(o)"
$ws.Range("I650").Value = "Primera"
$ws.Range("J650").Value = 120
$ws.Range("K650").Value = 6000
$ws.Range("L650").Value = 6500
$ws.Range("M650").Value = 6250
$ws.Range("N650").Value = "$/caja 10 unidades"
$ws.Range("O650").Value = "Provincia de Diguillín"
$ws.Range("P650").Value = 625
$ws.Range("Q650").Value = 10
$ws.Range("R650").Value = "Hortaliza"

# --- New row 651: Lechuga / Escarola ---
$ws.Range("A651").Value = 7
$ws.Range("B651").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C651").Value = "Ñuble"
$ws.Range("D651").Value = 44753
$ws.Range("E651").Value = 16
$ws.Range("F651").Value = 100112033
$ws.Range("G651").Value = "Lechuga"
$ws.Range("H651").Value = "Escarola"
$ws.Range("I651").Value = "Primera"
$ws.Range("J651").Value = 120
$ws.Range("K651").Value = 10000
$ws.Range("L651").Value = 11000
$ws.Range("M651").Value = 10500
$ws.Range("N651").Value = "$/caja 15 unidades"
$ws.Range("O651").Value = "Provincia del Elquí"
$ws.Range("P651").Value = 700
$ws.Range("Q651").Value = 15
$ws.Range("R651").Value = "Hortaliza"
